$wb = $excel.ActiveWorkbook

# Helper: set a cell's value as a literal TEXT string (not auto-converted to a
# number/date by Excel's type inference), without leaving any NumberFormat /
# style residue behind. We do this by entering a text-formula ("=""...""")
# and then collapsing the formula down to its static (text) result.
function Set-TextValue {
    param(
        $Sheet,
        [string]$Addr,
        [string]$Text
    )
    $escaped = $Text -replace '"', '""'
    $cell = $Sheet.Range($Addr)
    $cell.Formula = '="' + $escaped + '"'
    $cell.Value = $cell.Value()
}

# ---- Restricciones_del_follower ---------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $ws "A2" "2.6252662319059663 - x - 0.12098907152560284y"
Set-TextValue $ws "B2" "-0.6252662319059663"
Set-TextValue $ws "D2" "0.62"
Set-TextValue $ws "E2" "-2.0"
Set-TextValue $ws "F2" "-0.2"

Set-TextValue $ws "A3" "-22.11261681242672 - 0.25x + 4.597584717972911y"
Set-TextValue $ws "B3" "20.11261681242672"
Set-TextValue $ws "D3" "0.96"
Set-TextValue $ws "E3" "7.1"
Set-TextValue $ws "F3" "7.6"

Set-TextValue $ws "A4" "-0.423517021032052 + x - 0.32652093068454247y"
Set-TextValue $ws "B4" "-7.576482978967948"
Set-TextValue $ws "D4" "0.88"
Set-TextValue $ws "E4" "-1.1"
Set-TextValue $ws "F4" "-0.3"

Set-TextValue $ws "A5" "-2.8600000000000003 + x"
Set-TextValue $ws "B5" "0.03000000000000025"
Set-TextValue $ws "D5" "0.29"
Set-TextValue $ws "E5" "3.7"
Set-TextValue $ws "F5" "0.0"

Set-TextValue $ws "A6" "-26.579347319815504 + 5.38401368288933y"
Set-TextValue $ws "B6" "26.489347319815504"
Set-TextValue $ws "D6" "0.34"
Set-TextValue $ws "E6" "0"
Set-TextValue $ws "F6" "8.9"

# ---- Punto_modificado ----------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws "A2" "2.0300000000000002"
Set-TextValue $ws "B2" "4.92"

# ---- Vector_bf -------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $ws "A2" "-4.881894338088095"

# ---- Vector_BF -------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_BF")
Set-TextValue $ws "A2" "-3.825"
Set-TextValue $ws "A3" "-34.244002664411866"

# ---- Vector_Alpha (plain numeric cell, not text) ---------------------------
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 1.653041861369085
